# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the zh-cn and de-de handback rows, as part of a freshly
# generated handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 11:06:11"
$wsZhCn.Range("H2").Value = "2016-03-13 11:06:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 11:06:15"
$wsDeDe.Range("H2").Value = "2016-03-13 11:06:35"
